# fix(quickplay): it now load question excel by custom mode for quickplay
#
# Row 44 (the "quickplay" row) was only partially filled in (question,
# mode, operands) — it was missing the difficulty/equation/time columns
# that every other question row has. This fills those in and cleans up
# the stray border formatting that had been left on C44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C44 previously had a boxed/bordered style left over from editing;
# the other cells in the row (and the rest of the sheet) use the plain,
# unbordered style, so strip the border back off.
$ws.Cells.Item(44, 3).Borders.LineStyle = -4142   # xlLineStyleNone

# Fill in the rest of the "quickplay" question row: difficulty, equation, time.
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(44, 5).Value = "{a}+{b}"
$ws.Cells.Item(44, 6).Value = 30

# Row height had been left at a stray custom value (13.8) from the old
# formatting; auto-fit it back to the sheet's normal row height.
$ws.Rows.Item(44).AutoFit()

# Move the saved cursor/selection back up near the top of the sheet.
$ws.Range("B25").Select() | Out-Null
